$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new value in C3 (new row) as described by the diff
$ws.Range("C3").Value = "AzureBastiuoin"

# Update selection to match the post-edit state (C4, the cell below the new entry)
$ws.Range("C4").Select()
